# Apply the changes described in the commit:
#  - Add "Extra Duty Allowance" column to Sheet2 (column J), defaulting to 0,
#    with 100 for "Test Employee 3" (row 4)
#  - Add "Employee Name" and "Extra Duty Allowance" columns to the "Test Name"
#    sheet, referencing the chosen test record's employee name / amount, and
#    change the test id in column A
#  - Update the active selections / view state left behind by the edit session

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet2: add "Extra Duty Allowance" column (J)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Cells.Item(1, 10).Value = "Extra Duty Allowance"
$ws2.Cells.Item(1, 10).Font.Bold = $true

$ws2.Cells.Item(2, 10).Value = 0
$ws2.Cells.Item(3, 10).Value = 0
$ws2.Cells.Item(4, 10).Value = 100
$ws2.Cells.Item(5, 10).Value = 0
$ws2.Cells.Item(6, 10).Value = 0
$ws2.Cells.Item(7, 10).Value = 0
$ws2.Cells.Item(8, 10).Value = 0
$ws2.Cells.Item(9, 10).Value = 0
$ws2.Cells.Item(10, 10).Value = 0
$ws2.Cells.Item(11, 10).Value = 0
$ws2.Cells.Item(12, 10).Value = 0

$ws2.Columns.Item(10).ColumnWidth = 21.166666666666668

# ---------------------------------------------------------------------------
# Test Name sheet: choose record by test id, add employee name / extra duty
# allowance columns
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Test Name")

$ws3.Cells.Item(1, 2).Value = "Employee Name"
$ws3.Cells.Item(1, 2).Font.Bold = $true

$ws3.Cells.Item(2, 1).Value = "TkmwGCjQ"
$ws3.Cells.Item(2, 2).Value = "Test Employee 3"

$ws3.Cells.Item(1, 3).Value = "Extra Duty Allowance"
$ws3.Cells.Item(1, 3).Font.Bold = $true
$ws3.Cells.Item(2, 3).Value = "'100.00"

$ws3.Columns.Item(2).ColumnWidth = 21.736979166666668
$ws3.Columns.Item(3).ColumnWidth = 22.166666666666668

# ---------------------------------------------------------------------------
# View state: active cell / zoom left on each sheet
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A4").Select() | Out-Null

$ws2.Activate()
$ws2.Range("K16").Select() | Out-Null

$ws3.Activate()
$ws3.Range("E12").Select() | Out-Null
$excel.ActiveWindow.Zoom = 106
